$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ratings (column B) for the four new form responses
$ws.Range("B3").Value = 5
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 5
$ws.Range("B6").Value = 7

# Explain why (column C)
$ws.Range("C3").Value = "Idk"
$ws.Range("C4").Value = "It was bad"
$ws.Range("C5").Value = "Decent teacher"
$ws.Range("C6").Value = "Cool class"

# Extended Why (column D) and the first "would take again" answer (column E)
$ws.Range("D3").Value = "N/A"
$ws.Range("D4").Value = "Too much work"
$ws.Range("E3").Value = "Yes"
$ws.Range("D5").Value = "Taught well"
$ws.Range("D6").Value = "Went smoothly"

# Remaining "Would you take the class again" answers (column E)
$ws.Range("E4").Value = "No"
$ws.Range("E5").Value = "No"
$ws.Range("E6").Value = "Yes"
